$wb = $excel.ActiveWorkbook

# The "Descriptions" sheet gains a new 5th column (E1) with a question about
# the introductory QCM, and becomes the active/selected sheet.
$wsDescriptions = $wb.Worksheets.Item("Descriptions")
$wsFeatures = $wb.Worksheets.Item("Features")

$wsDescriptions.Range("E1").Value = "Question affichée dans la partie de QCM introductif"

# "Features" was the active tab before; now it's just a regular sheet whose
# selection moved to C1.
$wsFeatures.Range("C1").Select()

# Make "Descriptions" the active sheet/tab, with E1 selected.
$wsDescriptions.Activate()
$wsDescriptions.Range("E1").Select()
